$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) text renames ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C numeric data updates (rows 2-29) ---
$ws.Range("C2").Value = 9502.243585046588
$ws.Range("C3").Value = 1909.084588129339
$ws.Range("C4").Value = 4547.50930098406
$ws.Range("C5").Value = 4729.735976516416
$ws.Range("C6").Value = 10385.96443195552
$ws.Range("C7").Value = 1955.461557360978
$ws.Range("C8").Value = 4633.590358399045
$ws.Range("C9").Value = 5082.354756663512
$ws.Range("C10").Value = 10883.31535948899
$ws.Range("C11").Value = 2024.117324382548
$ws.Range("C12").Value = 13455.83781255333
$ws.Range("C13").Value = 4921.848409120176
$ws.Range("C14").Value = 5360.226632400601
$ws.Range("C15").Value = 2094.024217383061
$ws.Range("C16").Value = 5122.180090208862
$ws.Range("C17").Value = 5642.578115155247
$ws.Range("C18").Value = 2201.396847776877
$ws.Range("C19").Value = 5295.682695961288
$ws.Range("C20").Value = 5919.20956823756
$ws.Range("C21").Value = 2286.013198234259
$ws.Range("C22").Value = 5412.131646018807
$ws.Range("C23").Value = 5996.49696468919
$ws.Range("C24").Value = 2361.056581219794
$ws.Range("C25").Value = 5330.539154475424
$ws.Range("C26").Value = 6114.227214287786
$ws.Range("C27").Value = 2425.561644739583
$ws.Range("C28").Value = 5176.058803160127
$ws.Range("C29").Value = 6262.368904654469
